$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Minor precision correction of an existing value ---
$ws.Range("D55").Value2 = 0.112114197530864

# --- Insert two new rows (56 and 57) by copying the formatting of the
# --- previous two rows (54 and 55), which already carry the correct
# --- styles (s=6 on column D, s=7 on C/E/F/G/H, s=1 on A).
$ws.Rows("54:55").Copy()
$ws.Rows("56:57").Insert(-4121)

# --- Row 56: 2025-09-28, 四方坪站 ---
$ws.Range("A56").Value2 = 45928
$ws.Range("B56").Value2 = $ws.Range("B54").Value2
$ws.Range("C56").Formula = "=15788/127"
$ws.Range("D56").Formula = "=C56/(24*60)"
$ws.Range("E56").Formula = "=8093.57/127"
$ws.Range("F56").Formula = "=3124.72/127"
$ws.Range("G56").Formula = "=9093.57/(15788/60)"
$ws.Range("H56").Formula = "=386/127"

# --- Row 57: 2025-09-28, 高岭站 ---
$ws.Range("A57").Value2 = 45928
$ws.Range("B57").Value2 = $ws.Range("B55").Value2
$ws.Range("C57").Formula = "=5876/36"
$ws.Range("D57").Formula = "=C57/(24*60)"
$ws.Range("E57").Formula = "=5876.99/36"
$ws.Range("F57").Formula = "=1432.67/36"
$ws.Range("G57").Formula = "=5876.99/(7618/60)"
$ws.Range("H57").Formula = "=206/36"

# --- Rows 58-63: blank filler rows, formatted (C:H) like column C of row 54 ---
$ws.Range("C54").Copy()
$ws.Range("C58:H63").PasteSpecial(-4122)

# --- Update view state: scroll/selection moved down ---
$ws.Range("E68").Select() | Out-Null
